$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.705.22"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "2.802.87"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'355.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "'109.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  -1.24%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +5.39%  "

$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").Value = "'19.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "'7.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").Value = "3.243.13"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "2.805.44"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "51.676.27"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").Value = "'7.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.58%  "

$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").Value = "'13.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.52%  "

$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D24").Value = "'268.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("D25").Value = "'2.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("D28").Value = "'0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").Value = "'10.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("D30").Value = "'37.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.84%  "

$ws.Range("E31").Value = "  +2.58%  "

$ws.Range("D32").Value = "'6.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D34").Value = "'5.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.21%  "

$ws.Range("D35").Value = "'0.0447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.53%  "

$ws.Range("D36").Value = "'0.0859"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("D38").Value = "'18.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("D40").Value = "'3.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("E42").Value = "  -4.90%  "

$ws.Range("D43").Value = "'119.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "2.127.48"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("E47").Value = "  +6.37%  "

$ws.Range("D48").Value = "'3.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.45%  "

$ws.Range("D49").Value = "'0.907"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.29%  "

$ws.Range("D50").Value = "'5.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.76%  "

$ws.Range("E51").Value = "  +6.59%  "
